{"js": "const body = context.document.body;\n\nasync function replaceOne(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait replaceOne(\"2025-12-29 Monday\", \"2025-12-30 Tuesday\");\nawait replaceOne(\"350\u00f74=87, 2\", \"750\u00f79=83, 3\");\nawait replaceOne(\"453\u00f72=226, 1\", \"825\u00f75=165, 0\");\nawait replaceOne(\"418\u00f79=46, 4\", \"173\u00f77=24, 5\");\nawait replaceOne(\"848\u00f77=121, 1\", \"227\u00f74=56, 3\");\nawait replaceOne(\"535\u00f74=133, 3\", \"646\u00f78=80, 6\");\nawait replaceOne(\"993\u00f72=496, 1\", \"435\u00f78=54, 3\");\nawait replaceOne(\"663\u00f72=331, 1\", \"272\u00f72=136, 0\");\nawait replaceOne(\"183\u00f78=22, 7\", \"230\u00f73=76, 2\");\nawait replaceOne(\"732\u00f76=122, 0\", \"464\u00f78=58, 0\");\nawait replaceOne(\"760\u00f77=108, 4\", \"718\u00f72=359, 0\");\nawait replaceOne(\"762\u00f72=381, 0\", \"407\u00f73=135, 2\");\nawait replaceOne(\"176\u00f73=58, 2\", \"315\u00f78=39, 3\");\nawait replaceOne(\"903\u00f74=225, 3\", \"603\u00f76=100, 3\");\nawait replaceOne(\"944\u00f76=157, 2\", \"730\u00f78=91, 2\");\nawait replaceOne(\"371\u00f74=92, 3\", \"237\u00f75=47, 2\");\nawait replaceOne(\"136\u00f78=17, 0\", \"838\u00f79=93, 1\");\nawait replaceOne(\"809\u00f78=101, 1\", \"814\u00f77=116, 2\");\nawait replaceOne(\"761\u00f72=380, 1\", \"552\u00f75=110, 2\");\nawait replaceOne(\"107\u00f75=21, 2\", \"159\u00f78=19, 7\");\nawait replaceOne(\"264\u00f74=66, 0\", \"906\u00f77=129, 3\");\nawait replaceOne(\"769\u00f72=384, 1\", \"120\u00f76=20, 0\");\nawait replaceOne(\"244\u00f78=30, 4\", \"460\u00f72=230, 0\");\nawait replaceOne(\"686\u00f73=228, 2\", \"817\u00f75=163, 2\");\nawait replaceOne(\"836\u00f75=167, 1\", \"830\u00f76=138, 2\");\nawait replaceOne(\"488\u00f73=162, 2\", \"273\u00f79=30, 3\");\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n  $range = $d.Content\n  $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\nReplace-Text '2025-12-29 Monday' '2025-12-30 Tuesday'\nReplace-Text '350\u00f74=87, 2' '750\u00f79=83, 3'\nReplace-Text '453\u00f72=226, 1' '825\u00f75=165, 0'\nReplace-Text '418\u00f79=46, 4' '173\u00f77=24, 5'\nReplace-Text '848\u00f77=121, 1' '227\u00f74=56, 3'\nReplace-Text '535\u00f74=133, 3' '646\u00f78=80, 6'\nReplace-Text '993\u00f72=496, 1' '435\u00f78=54, 3'\nReplace-Text '663\u00f72=331, 1' '272\u00f72=136, 0'\nReplace-Text '183\u00f78=22, 7' '230\u00f73=76, 2'\nReplace-Text '732\u00f76=122, 0' '464\u00f78=58, 0'\nReplace-Text '760\u00f77=108, 4' '718\u00f72=359, 0'\nReplace-Text '762\u00f72=381, 0' '407\u00f73=135, 2'\nReplace-Text '176\u00f73=58, 2' '315\u00f78=39, 3'\nReplace-Text '903\u00f74=225, 3' '603\u00f76=100, 3'\nReplace-Text '944\u00f76=157, 2' '730\u00f78=91, 2'\nReplace-Text '371\u00f74=92, 3' '237\u00f75=47, 2'\nReplace-Text '136\u00f78=17, 0' '838\u00f79=93, 1'\nReplace-Text '809\u00f78=101, 1' '814\u00f77=116, 2'\nReplace-Text '761\u00f72=380, 1' '552\u00f75=110, 2'\nReplace-Text '107\u00f75=21, 2' '159\u00f78=19, 7'\nReplace-Text '264\u00f74=66, 0' '906\u00f77=129, 3'\nReplace-Text '769\u00f72=384, 1' '120\u00f76=20, 0'\nReplace-Text '244\u00f78=30, 4' '460\u00f72=230, 0'\nReplace-Text '686\u00f73=228, 2' '817\u00f75=163, 2'\nReplace-Text '836\u00f75=167, 1' '830\u00f76=138, 2'\nReplace-Text '488\u00f73=162, 2' '273\u00f79=30, 3'\n"}
